$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("UniSolder52C_front")

# R70 must actually be 0 ohm, so it no longer belongs with the rest of the
# 1.8k resistors ("R23, R27, R28, R33, R44, R57, R70"). Pull it out of that
# combined row and give it its own BOM line.

# 1) Shrink the existing row 17 designator list / quantity (R70 removed,
#    leaving the other 1.8k resistors).
$ws.Range("A17").Value = "R23, R27, R28, R33, R44, R57"
$ws.Range("B17").Value = 6

# Re-assigning a text value that doesn't look like a number drops the
# "quote prefix" formatting Excel had applied to this (already-text) cell;
# restore it (and the rest of the cell's look) by pasting the format from a
# neighbouring text cell in the same row.
$ws.Range("C17").Copy()
$ws.Range("A17").PasteSpecial(-4122)

# 2) Insert a new row for R70 right before the old row 28 ("Rc2"),
#    pushing that row (and everything below it) down by one.
$ws.Rows.Item(28).Insert()

# Keep "0" (and "0805") as text, matching how this BOM stores every other
# Value/Package entry, rather than letting Excel coerce them to numbers.
$ws.Range("C28").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("G28").NumberFormat = "@"

$ws.Range("A28").Value = "R70"
$ws.Range("B28").Value = 1
$ws.Range("C28").Value = "0"
$ws.Range("D28").Value = "0"
$ws.Range("E28").Value = "Resistor"
$ws.Range("F28").Value = "0805 Resistor"
$ws.Range("G28").Value = "0805"

# Match the look (borders/shading/alignment/quote-prefix) of the other
# resistor rows, e.g. the row just above.
$ws.Range("A27:G27").Copy()
$ws.Range("A28:G28").PasteSpecial(-4122)
